# [FIX] update data formatting for consistency
#
# This script reproduces, via Excel COM automation, the edits that were made
# to financial_data/planet.xlsx:
#   1. Relabel five header cells on the "consolidated" sheet so their shared
#      strings become distinct from the generic labels used on other sheets
#      ("Total assets" -> "Total Assets", etc.)
#   2. Remove the oldest quarter (row 4, dated 2022-02-28 / 44592) from the
#      "consolidated" sheet, shifting the remaining quarters up by one row.
#   3. Update the selected cell / active sheet bookmarks that Excel stores
#      per worksheet, ending with the "income" sheet active.
#   4. Normalize the formatting (remove a redundant "no-fill" flag) of the
#      row-label cells in column A of the "pre-income" sheet.

$wb = $excel.ActiveWorkbook

$wsConsolidated   = $wb.Worksheets.Item(1)   # "consolidated"
$wsIncome         = $wb.Worksheets.Item(2)   # "income"
$wsPreConsolidated = $wb.Worksheets.Item(3)  # "pre-consolidation"
$wsPreIncome      = $wb.Worksheets.Item(4)   # "pre-income"

# ---------------------------------------------------------------------
# 1. Relabel the summary-total header cells on the "consolidated" sheet.
#    Setting them in this order makes the newly created shared strings
#    land in the same order they appear in the saved workbook:
#    Total Assets, Total Liabilities, Total Current Assets,
#    Total Current Liabilities, Total Equity.
# ---------------------------------------------------------------------
$wsConsolidated.Range("N1").Value = "Total Assets"
$wsConsolidated.Range("AB1").Value = "Total Liabilities"
$wsConsolidated.Range("F1").Value = "Total Current Assets"
$wsConsolidated.Range("T1").Value = "Total Current Liabilities"
$wsConsolidated.Range("AF1").Value = "Total Equity"

# ---------------------------------------------------------------------
# 2. Delete the oldest quarter's row; the rows below shift up, keeping
#    their original per-cell styles and values intact.
# ---------------------------------------------------------------------
$wsConsolidated.Rows(4).Delete()

# ---------------------------------------------------------------------
# 3. Update remembered selections / active sheet.
#    Order matters: the last sheet Activate()'d becomes the workbook's
#    active tab (tabSelected="1" / bookViews activeTab).
# ---------------------------------------------------------------------
$wsConsolidated.Activate()
$wsConsolidated.Range("A6").Select()

$wsIncome.Activate()
$wsIncome.Range("E22").Select()

# ---------------------------------------------------------------------
# 4. Normalize column-A label formatting on the "pre-income" sheet by
#    clearing the (already invisible, no-op) explicit "no fill" flag
#    that was left over on these cells, collapsing them onto the same
#    style already used by equivalent rows (e.g. A6, A11, A12).
# ---------------------------------------------------------------------
$preIncomeLabelCells = @("A2", "A3", "A4", "A5", "A7", "A8", "A9", "A10", "A13", "A14", "A15", "A16", "A17", "A18")
foreach ($cellAddr in $preIncomeLabelCells) {
    $wsPreIncome.Range($cellAddr).Interior.Pattern = -4142   # xlNone
}
